$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text run formatting) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "49"

$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 10).Text = "12/4/2023"
$weekCell.Characters(47, 9).Text = "12/10/2023"

# --- Crime-statistics grid updates ---
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = '#,##0'
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 33.333333333333
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 33.333333333333
$ws.Range("N14").Value = -76.470588235294
$ws.Range("M15").Value = -23.333333333333
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 7.407407407407
$ws.Range("I16").Value = 338
$ws.Range("J16").Value = 296
$ws.Range("K16").Value = 14.189189189189
$ws.Range("L16").Value = 43.829787234042
$ws.Range("M16").Value = 3.363914373088
$ws.Range("N16").Value = -72.675828617623
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 40.74074074074
$ws.Range("I17").Value = 455
$ws.Range("J17").Value = 391
$ws.Range("K17").Value = 16.368286445012
$ws.Range("L17").Value = 22.972972972973
$ws.Range("M17").Value = 47.727272727272
$ws.Range("N17").Value = 11.519607843137
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 132
$ws.Range("J18").Value = 141
$ws.Range("K18").Value = -6.382978723404
$ws.Range("L18").Value = -4.347826086956
$ws.Range("M18").Value = -52.51798561151
$ws.Range("N18").Value = -92.967501331912
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 16
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -12.328767123287
$ws.Range("I19").Value = 808
$ws.Range("J19").Value = 925
$ws.Range("K19").Value = -12.648648648648
$ws.Range("L19").Value = 9.189189189189
$ws.Range("M19").Value = 67.983367983368
$ws.Range("N19").Value = -40.588235294117
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 326
$ws.Range("J20").Value = 302
$ws.Range("K20").Value = 7.947019867549
$ws.Range("L20").Value = 68.911917098445
$ws.Range("M20").Value = 41.125541125541
$ws.Range("N20").Value = -84.651600753295
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 169
$ws.Range("G21").Value = 167
$ws.Range("H21").Value = 1.197604790419
$ws.Range("I21").Value = 2086
$ws.Range("J21").Value = 2094
$ws.Range("K21").Value = -0.382043935052
$ws.Range("L21").Value = 21.917007597896
$ws.Range("M21").Value = 25.814234016887
$ws.Range("N21").Value = -70.436507936507
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = -50
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G22").Value = 2
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("H22").Value = 250
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I22").Value = 84
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 40
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 200
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = -19.565217391304
$ws.Range("G24").Value = 192
$ws.Range("H24").Value = -11.979166666666
$ws.Range("I24").Value = 1890
$ws.Range("J24").Value = 1856
$ws.Range("K24").Value = 1.831896551724
$ws.Range("L24").Value = 44.495412844036
$ws.Range("M24").Value = 75.487465181058
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = 52.631578947368
$ws.Range("F25").Value = 89
$ws.Range("G25").Value = 74
$ws.Range("H25").Value = 20.27027027027
$ws.Range("I25").Value = 939
$ws.Range("J25").Value = 862
$ws.Range("K25").Value = 8.932714617169
$ws.Range("L25").Value = 17.228464419475
$ws.Range("M25").Value = 7.560137457044
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 49
$ws.Range("K26").Value = -8.163265306122
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 125
$ws.Range("I27").Value = 138
$ws.Range("J27").Value = 106
$ws.Range("K27").Value = 30.188679245283
$ws.Range("L27").Value = 38
$ws.Range("N28").Value = -91.525423728813
$ws.Range("N29").Value = -92.727272727272
